$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 692-701: vocabulary words (dans la rue lesson)
$ws.Cells.Item(692, 2).Value = "on"
$ws.Cells.Item(692, 3).Value = "dans"
$ws.Cells.Item(692, 4).Value = "dɑ̃"

$ws.Cells.Item(693, 2).Value = "the road"
$ws.Cells.Item(693, 3).Value = "la rue"
$ws.Cells.Item(693, 4).Value = "la ʁy"

$ws.Cells.Item(694, 2).Value = "some"
$ws.Cells.Item(694, 3).Value = "du"
$ws.Cells.Item(694, 4).Value = "dy"

$ws.Cells.Item(695, 2).Value = "fire"
$ws.Cells.Item(695, 3).Value = "feu"
$ws.Cells.Item(695, 4).Value = "fø"

$ws.Cells.Item(696, 2).Value = "smoke"
$ws.Cells.Item(696, 3).Value = "fume"
$ws.Cells.Item(696, 4).Value = "fym"

$ws.Cells.Item(697, 2).Value = "plus/more"
$ws.Cells.Item(697, 3).Value = "plus"
$ws.Cells.Item(697, 4).Value = "ply"

$ws.Cells.Item(698, 2).Value = "table"
$ws.Cells.Item(698, 3).Value = "la table"
$ws.Cells.Item(698, 4).Value = "tabl"
$ws.Cells.Item(698, 5).Value = "F"

$ws.Cells.Item(699, 2).Value = "route"
$ws.Cells.Item(699, 3).Value = "la route"
$ws.Cells.Item(699, 4).Value = "ʁut"
$ws.Cells.Item(699, 5).Value = "F"

$ws.Cells.Item(700, 2).Value = "tree"
$ws.Cells.Item(700, 3).Value = "l'arbre"
$ws.Cells.Item(700, 4).Value = "aʁbʁ"
$ws.Cells.Item(700, 5).Value = "M"

$ws.Cells.Item(701, 2).Value = "shop"
$ws.Cells.Item(701, 3).Value = "magasin"
$ws.Cells.Item(701, 4).Value = "ma.ɡa.zɛ̃"
$ws.Cells.Item(701, 5).Value = "M"

# Rows 702-721: numbers 1-20
$nCell = $ws.Cells.Item(702, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 1
$ws.Cells.Item(702, 3).Value = "un"
$ws.Cells.Item(702, 4).Value = "œ̃ / ɛ̃"

$nCell = $ws.Cells.Item(703, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 2
$ws.Cells.Item(703, 3).Value = "deux"
$ws.Cells.Item(703, 4).Value = "dø"

$nCell = $ws.Cells.Item(704, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 3
$ws.Cells.Item(704, 3).Value = "trois"
$ws.Cells.Item(704, 4).Value = "tʁwɑ"

$nCell = $ws.Cells.Item(705, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 4
$ws.Cells.Item(705, 3).Value = "quatre"
$ws.Cells.Item(705, 4).Value = "katʁ"

$nCell = $ws.Cells.Item(706, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 5
$ws.Cells.Item(706, 3).Value = "cinq"
$ws.Cells.Item(706, 4).Value = "sɛ̃k"

$nCell = $ws.Cells.Item(707, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 6
$ws.Cells.Item(707, 3).Value = "six"
$ws.Cells.Item(707, 4).Value = "sis"

$nCell = $ws.Cells.Item(708, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 7
$ws.Cells.Item(708, 3).Value = "sept"
$ws.Cells.Item(708, 4).Value = "sɛt"

$nCell = $ws.Cells.Item(709, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 8
$ws.Cells.Item(709, 3).Value = "huit"
$ws.Cells.Item(709, 4).Value = "ɥit"

$nCell = $ws.Cells.Item(710, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 9
$ws.Cells.Item(710, 3).Value = "neuf"
$ws.Cells.Item(710, 4).Value = "nœf"

$nCell = $ws.Cells.Item(711, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 10
$ws.Cells.Item(711, 3).Value = "dix"
$ws.Cells.Item(711, 4).Value = "dis"

$nCell = $ws.Cells.Item(712, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 11
$ws.Cells.Item(712, 3).Value = "onze"
$ws.Cells.Item(712, 4).Value = "ɔ̃z"

$nCell = $ws.Cells.Item(713, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 12
$ws.Cells.Item(713, 3).Value = "douze"
$ws.Cells.Item(713, 4).Value = "duz"

$nCell = $ws.Cells.Item(714, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 13
$ws.Cells.Item(714, 3).Value = "treize"
$ws.Cells.Item(714, 4).Value = "tʁɛz"

$nCell = $ws.Cells.Item(715, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 14
$ws.Cells.Item(715, 3).Value = "quatorze"
$ws.Cells.Item(715, 4).Value = "ka.tɔʁz"

$nCell = $ws.Cells.Item(716, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 15
$ws.Cells.Item(716, 3).Value = "quinze"
$ws.Cells.Item(716, 4).Value = "kɛ̃z"

$nCell = $ws.Cells.Item(717, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 16
$ws.Cells.Item(717, 3).Value = "seize"
$ws.Cells.Item(717, 4).Value = "sɛz"

$nCell = $ws.Cells.Item(718, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 17
$ws.Cells.Item(718, 3).Value = "dix-sept"

$nCell = $ws.Cells.Item(719, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 18
$ws.Cells.Item(719, 3).Value = "dix-huit"
$ws.Cells.Item(719, 4).Value = "di.zɥit"

$nCell = $ws.Cells.Item(720, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 19
$ws.Cells.Item(720, 3).Value = "dix-neuf"

$nCell = $ws.Cells.Item(721, 2)
$nCell.Style = "Normal"
$nCell.Font.Size = 20
$nCell.Value = 20
$ws.Cells.Item(721, 3).Value = "vignt"
$ws.Cells.Item(721, 4).Value = "vɛ̃"
